$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Selplg"
$ws.Cells.Item(2, 3).Value = "Sell"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 3.137719
$ws.Cells.Item(2, 8).Value = 9.413157
$ws.Cells.Item(2, 9).Value = 0.02008359063640063
$ws.Cells.Item(2, 10).Value = 0.02008359063640063
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.2056386666666667
$ws.Cells.Item(2, 14).Value = 0.616916
$ws.Cells.Item(2, 15).Value = 0.004186411275012692
$ws.Cells.Item(2, 16).Value = 0.004186411275012692
$ws.Cells.Item(2, 17).Value = 0.6452363515346667
$ws.Cells.Item(2, 18).Value = 5.807127163812
$ws.Cells.Item(2, 19).Value = 0.000084078170282966943529373571
$ws.Cells.Item(2, 20).Value = 0.000084078170282966929976846415

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Selplg"
$ws.Cells.Item(3, 3).Value = "Sell"
$ws.Cells.Item(3, 4).Value = "M1"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3.137719
$ws.Cells.Item(3, 8).Value = 9.413157
$ws.Cells.Item(3, 9).Value = 0.02008359063640063
$ws.Cells.Item(3, 10).Value = 0.02008359063640063
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.974950666666667
$ws.Cells.Item(3, 14).Value = 8.924852
$ws.Cells.Item(3, 15).Value = 0.06056432486857137
$ws.Cells.Item(3, 16).Value = 0.06056432486857137
$ws.Cells.Item(3, 17).Value = 9.334559230862666
$ws.Cells.Item(3, 18).Value = 84.011033077764
$ws.Cells.Item(3, 19).Value = 0.001216349107830366
$ws.Cells.Item(3, 20).Value = 0.001216349107830366

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Selplg"
$ws.Cells.Item(4, 3).Value = "Sell"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.137719
$ws.Cells.Item(4, 8).Value = 9.413157
$ws.Cells.Item(4, 9).Value = 0.02008359063640063
$ws.Cells.Item(4, 10).Value = 0.02008359063640063
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 45.939923
$ws.Cells.Item(4, 14).Value = 137.819769
$ws.Cells.Item(4, 15).Value = 0.935249263856416
$ws.Cells.Item(4, 16).Value = 0.935249263856416
$ws.Cells.Item(4, 17).Value = 144.146569255637
$ws.Cells.Item(4, 18).Value = 1297.319123300733
$ws.Cells.Item(4, 19).Value = 0.0187831633582873
$ws.Cells.Item(4, 20).Value = 0.0187831633582873

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Selplg"
$ws.Cells.Item(5, 3).Value = "Sell"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.516030666666667
$ws.Cells.Item(5, 8).Value = 7.548092
$ws.Cells.Item(5, 9).Value = 0.01610435158086607
$ws.Cells.Item(5, 10).Value = 0.01610435158086607
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.2056386666666667
$ws.Cells.Item(5, 14).Value = 0.616916
$ws.Cells.Item(5, 15).Value = 0.004186411275012692
$ws.Cells.Item(5, 16).Value = 0.004186411275012692
$ws.Cells.Item(5, 17).Value = 0.5173931915857778
$ws.Cells.Item(5, 18).Value = 4.656538724272
$ws.Cells.Item(5, 19).Value = 0.000067419439034906180653605734
$ws.Cells.Item(5, 20).Value = 0.000067419439034906180653605734

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Selplg"
$ws.Cells.Item(6, 3).Value = "Sell"
$ws.Cells.Item(6, 4).Value = "M1"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.516030666666667
$ws.Cells.Item(6, 8).Value = 7.548092
$ws.Cells.Item(6, 9).Value = 0.01610435158086607
$ws.Cells.Item(6, 10).Value = 0.01610435158086607
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.974950666666667
$ws.Cells.Item(6, 14).Value = 8.924852
$ws.Cells.Item(6, 15).Value = 0.06056432486857137
$ws.Cells.Item(6, 16).Value = 0.06056432486857137
$ws.Cells.Item(6, 17).Value = 7.485067109153777
$ws.Cells.Item(6, 18).Value = 67.36560398238399
$ws.Cells.Item(6, 19).Value = 0.0009753491809412634
$ws.Cells.Item(6, 20).Value = 0.0009753491809412634

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Selplg"
$ws.Cells.Item(7, 3).Value = "Sell"
$ws.Cells.Item(7, 4).Value = "M2"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.516030666666667
$ws.Cells.Item(7, 8).Value = 7.548092
$ws.Cells.Item(7, 9).Value = 0.01610435158086607
$ws.Cells.Item(7, 10).Value = 0.01610435158086607
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 45.939923
$ws.Cells.Item(7, 14).Value = 137.819769
$ws.Cells.Item(7, 15).Value = 0.935249263856416
$ws.Cells.Item(7, 16).Value = 0.935249263856416
$ws.Cells.Item(7, 17).Value = 115.5862550923053
$ws.Cells.Item(7, 18).Value = 1040.276295830748
$ws.Cells.Item(7, 19).Value = 0.0150615829608899
$ws.Cells.Item(7, 20).Value = 0.0150615829608899

# Row 8
$ws.Cells.Item(8, 1).Value = "M1"
$ws.Cells.Item(8, 2).Value = "Selplg"
$ws.Cells.Item(8, 3).Value = "Sell"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 72.66191099999999
$ws.Cells.Item(8, 8).Value = 217.985733
$ws.Cells.Item(8, 9).Value = 0.4650869231382975
$ws.Cells.Item(8, 10).Value = 0.4650869231382975
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.2056386666666667
$ws.Cells.Item(8, 14).Value = 0.616916
$ws.Cells.Item(8, 15).Value = 0.004186411275012692
$ws.Cells.Item(8, 16).Value = 0.004186411275012692
$ws.Cells.Item(8, 17).Value = 14.942098495492
$ws.Cells.Item(8, 18).Value = 134.478886459428
$ws.Cells.Item(8, 19).Value = 0.00194704513888713
$ws.Cells.Item(8, 20).Value = 0.00194704513888713

# Row 9
$ws.Cells.Item(9, 1).Value = "M1"
$ws.Cells.Item(9, 2).Value = "Selplg"
$ws.Cells.Item(9, 3).Value = "Sell"
$ws.Cells.Item(9, 4).Value = "M1"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 72.66191099999999
$ws.Cells.Item(9, 8).Value = 217.985733
$ws.Cells.Item(9, 9).Value = 0.4650869231382975
$ws.Cells.Item(9, 10).Value = 0.4650869231382975
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.974950666666667
$ws.Cells.Item(9, 14).Value = 8.924852
$ws.Cells.Item(9, 15).Value = 0.06056432486857137
$ws.Cells.Item(9, 16).Value = 0.06056432486857137
$ws.Cells.Item(9, 17).Value = 216.165600570724
$ws.Cells.Item(9, 18).Value = 1945.490405136516
$ws.Cells.Item(9, 19).Value = 0.02816767550507213
$ws.Cells.Item(9, 20).Value = 0.02816767550507213

# Row 10
$ws.Cells.Item(10, 1).Value = "M1"
$ws.Cells.Item(10, 2).Value = "Selplg"
$ws.Cells.Item(10, 3).Value = "Sell"
$ws.Cells.Item(10, 4).Value = "M2"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 72.66191099999999
$ws.Cells.Item(10, 8).Value = 217.985733
$ws.Cells.Item(10, 9).Value = 0.4650869231382975
$ws.Cells.Item(10, 10).Value = 0.4650869231382975
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 45.939923
$ws.Cells.Item(10, 14).Value = 137.819769
$ws.Cells.Item(10, 15).Value = 0.935249263856416
$ws.Cells.Item(10, 16).Value = 0.935249263856416
$ws.Cells.Item(10, 17).Value = 3338.082596372853
$ws.Cells.Item(10, 18).Value = 30042.74336735568
$ws.Cells.Item(10, 19).Value = 0.4349722024943383
$ws.Cells.Item(10, 20).Value = 0.4349722024943383

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Selplg"
$ws.Cells.Item(11, 3).Value = "Sell"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 76.94530466666667
$ws.Cells.Item(11, 8).Value = 230.835914
$ws.Cells.Item(11, 9).Value = 0.4925036309237572
$ws.Cells.Item(11, 10).Value = 0.4925036309237571
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.2056386666666667
$ws.Cells.Item(11, 14).Value = 0.616916
$ws.Cells.Item(11, 15).Value = 0.004186411275012692
$ws.Cells.Item(11, 16).Value = 0.004186411275012692
$ws.Cells.Item(11, 17).Value = 15.82292985791378
$ws.Cells.Item(11, 18).Value = 142.406368721224
$ws.Cells.Item(11, 19).Value = 0.002061822753483907
$ws.Cells.Item(11, 20).Value = 0.002061822753483906

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Selplg"
$ws.Cells.Item(12, 3).Value = "Sell"
$ws.Cells.Item(12, 4).Value = "M1"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 76.94530466666667
$ws.Cells.Item(12, 8).Value = 230.835914
$ws.Cells.Item(12, 9).Value = 0.4925036309237572
$ws.Cells.Item(12, 10).Value = 0.4925036309237571
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 2.974950666666667
$ws.Cells.Item(12, 14).Value = 8.924852
$ws.Cells.Item(12, 15).Value = 0.06056432486857137
$ws.Cells.Item(12, 16).Value = 0.06056432486857137
$ws.Cells.Item(12, 17).Value = 228.9084854149698
$ws.Cells.Item(12, 18).Value = 2060.176368734728
$ws.Cells.Item(12, 19).Value = 0.0298281499022174
$ws.Cells.Item(12, 20).Value = 0.0298281499022174

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Selplg"
$ws.Cells.Item(13, 3).Value = "Sell"
$ws.Cells.Item(13, 4).Value = "M2"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 76.94530466666667
$ws.Cells.Item(13, 8).Value = 230.835914
$ws.Cells.Item(13, 9).Value = 0.4925036309237572
$ws.Cells.Item(13, 10).Value = 0.4925036309237571
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 45.939923
$ws.Cells.Item(13, 14).Value = 137.819769
$ws.Cells.Item(13, 15).Value = 0.935249263856416
$ws.Cells.Item(13, 16).Value = 0.935249263856416
$ws.Cells.Item(13, 17).Value = 3534.861371598208
$ws.Cells.Item(13, 18).Value = 31813.75234438387
$ws.Cells.Item(13, 19).Value = 0.4606136582680559
$ws.Cells.Item(13, 20).Value = 0.4606136582680558

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Selplg"
$ws.Cells.Item(14, 3).Value = "Sell"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.972004
$ws.Cells.Item(14, 8).Value = 2.916012
$ws.Cells.Item(14, 9).Value = 0.006221503720678607
$ws.Cells.Item(14, 10).Value = 0.006221503720678606
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.2056386666666667
$ws.Cells.Item(14, 14).Value = 0.616916
$ws.Cells.Item(14, 15).Value = 0.004186411275012692
$ws.Cells.Item(14, 16).Value = 0.004186411275012692
$ws.Cells.Item(14, 17).Value = 0.1998816065546667
$ws.Cells.Item(14, 18).Value = 1.798934458992
$ws.Cells.Item(14, 19).Value = 0.000026045773323782338732083549
$ws.Cells.Item(14, 20).Value = 0.000026045773323782328567688182

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Selplg"
$ws.Cells.Item(15, 3).Value = "Sell"
$ws.Cells.Item(15, 4).Value = "M1"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.972004
$ws.Cells.Item(15, 8).Value = 2.916012
$ws.Cells.Item(15, 9).Value = 0.006221503720678607
$ws.Cells.Item(15, 10).Value = 0.006221503720678606
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 2.974950666666667
$ws.Cells.Item(15, 14).Value = 8.924852
$ws.Cells.Item(15, 15).Value = 0.06056432486857137
$ws.Cells.Item(15, 16).Value = 0.06056432486857137
$ws.Cells.Item(15, 17).Value = 2.891663947802666
$ws.Cells.Item(15, 18).Value = 26.024975530224
$ws.Cells.Item(15, 19).Value = 0.0003768011725102046
$ws.Cells.Item(15, 20).Value = 0.0003768011725102046

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Selplg"
$ws.Cells.Item(16, 3).Value = "Sell"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.972004
$ws.Cells.Item(16, 8).Value = 2.916012
$ws.Cells.Item(16, 9).Value = 0.006221503720678607
$ws.Cells.Item(16, 10).Value = 0.006221503720678606
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 45.939923
$ws.Cells.Item(16, 14).Value = 137.819769
$ws.Cells.Item(16, 15).Value = 0.935249263856416
$ws.Cells.Item(16, 16).Value = 0.935249263856416
$ws.Cells.Item(16, 17).Value = 44.653788915692
$ws.Cells.Item(16, 18).Value = 401.884100241228
$ws.Cells.Item(16, 19).Value = 0.005818656774844621
$ws.Cells.Item(16, 20).Value = 0.00581865677484462
